# Adds the "opposing palace" (cung đối) interpretations for palaces
# without a main star (vô chính diệu) - one row per single star and per
# pair of stars borrowed from the opposite palace, appended after the
# existing data block (with a blank separator row, matching the sheet's
# existing section-break convention).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStrings = @(
    "Tử Vi tọa thủ tại cung đối Điền Trạch",
    "Thiên Cơ tọa thủ tại cung đối Điền Trạch",
    "Thái Dương tọa thủ tại cung đối Điền Trạch",
    "Vũ Khúc tọa thủ tại cung đối Điền Trạch",
    "Thiên Đồng tọa thủ tại cung đối Điền Trạch",
    "Liêm Trinh tọa thủ tại cung đối Điền Trạch",
    "Thiên Phủ tọa thủ tại cung đối Điền Trạch",
    "Thái Âm tọa thủ tại cung đối Điền Trạch",
    "Tham Lang tọa thủ tại cung đối Điền Trạch",
    "Cự Môn tọa thủ tại cung đối Điền Trạch",
    "Thiên Tướng tọa thủ tại cung đối Điền Trạch",
    "Thiên Lương tọa thủ tại cung đối Điền Trạch",
    "Thất Sát tọa thủ tại cung đối Điền Trạch",
    "Phá Quân tọa thủ tại cung đối Điền Trạch",
    "Tử Vi đồng cung Thiên Cơ tại cung đối Điền Trạch",
    "Tử Vi đồng cung Thái Dương tại cung đối Điền Trạch",
    "Tử Vi đồng cung Vũ Khúc tại cung đối Điền Trạch",
    "Tử Vi đồng cung Thiên Đồng tại cung đối Điền Trạch",
    "Tử Vi đồng cung Liêm Trinh tại cung đối Điền Trạch",
    "Tử Vi đồng cung Thiên Phủ tại cung đối Điền Trạch",
    "Tử Vi đồng cung Thái Âm tại cung đối Điền Trạch",
    "Tử Vi đồng cung Tham Lang tại cung đối Điền Trạch",
    "Tử Vi đồng cung Cự Môn tại cung đối Điền Trạch",
    "Tử Vi đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Tử Vi đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Tử Vi đồng cung Thất Sát tại cung đối Điền Trạch",
    "Tử Vi đồng cung Phá Quân tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Thái Dương tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Vũ Khúc tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Thiên Đồng tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Liêm Trinh tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Thiên Phủ tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Thái Âm tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Tham Lang tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Cự Môn tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Thất Sát tại cung đối Điền Trạch",
    "Thiên Cơ đồng cung Phá Quân tại cung đối Điền Trạch",
    "Thái Dương đồng cung Vũ Khúc tại cung đối Điền Trạch",
    "Thái Dương đồng cung Thiên Đồng tại cung đối Điền Trạch",
    "Thái Dương đồng cung Liêm Trinh tại cung đối Điền Trạch",
    "Thái Dương đồng cung Thiên Phủ tại cung đối Điền Trạch",
    "Thái Dương đồng cung Thái Âm tại cung đối Điền Trạch",
    "Thái Dương đồng cung Tham Lang tại cung đối Điền Trạch",
    "Thái Dương đồng cung Cự Môn tại cung đối Điền Trạch",
    "Thái Dương đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Thái Dương đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Thái Dương đồng cung Thất Sát tại cung đối Điền Trạch",
    "Thái Dương đồng cung Phá Quân tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Thiên Đồng tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Liêm Trinh tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Thiên Phủ tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Thái Âm tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Tham Lang tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Cự Môn tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Thất Sát tại cung đối Điền Trạch",
    "Vũ Khúc đồng cung Phá Quân tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Liêm Trinh tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Thiên Phủ tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Thái Âm tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Tham Lang tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Cự Môn tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Thất Sát tại cung đối Điền Trạch",
    "Thiên Đồng đồng cung Phá Quân tại cung đối Điền Trạch",
    "Liêm Trinh đồng cung Thiên Phủ tại cung đối Điền Trạch",
    "Liêm Trinh đồng cung Thái Âm tại cung đối Điền Trạch",
    "Liêm Trinh đồng cung Tham Lang tại cung đối Điền Trạch",
    "Liêm Trinh đồng cung Cự Môn tại cung đối Điền Trạch",
    "Liêm Trinh đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Liêm Trinh đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Liêm Trinh đồng cung Thất Sát tại cung đối Điền Trạch",
    "Liêm Trinh đồng cung Phá Quân tại cung đối Điền Trạch",
    "Thiên Phủ đồng cung Thái Âm tại cung đối Điền Trạch",
    "Thiên Phủ đồng cung Tham Lang tại cung đối Điền Trạch",
    "Thiên Phủ đồng cung Cự Môn tại cung đối Điền Trạch",
    "Thiên Phủ đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Thiên Phủ đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Thiên Phủ đồng cung Thất Sát tại cung đối Điền Trạch",
    "Thiên Phủ đồng cung Phá Quân tại cung đối Điền Trạch",
    "Thái Âm đồng cung Tham Lang tại cung đối Điền Trạch",
    "Thái Âm đồng cung Cự Môn tại cung đối Điền Trạch",
    "Thái Âm đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Thái Âm đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Thái Âm đồng cung Thất Sát tại cung đối Điền Trạch",
    "Thái Âm đồng cung Phá Quân tại cung đối Điền Trạch",
    "Tham Lang đồng cung Cự Môn tại cung đối Điền Trạch",
    "Tham Lang đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Tham Lang đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Tham Lang đồng cung Thất Sát tại cung đối Điền Trạch",
    "Tham Lang đồng cung Phá Quân tại cung đối Điền Trạch",
    "Cự Môn đồng cung Thiên Tướng tại cung đối Điền Trạch",
    "Cự Môn đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Cự Môn đồng cung Thất Sát tại cung đối Điền Trạch",
    "Cự Môn đồng cung Phá Quân tại cung đối Điền Trạch",
    "Thiên Tướng đồng cung Thiên Lương tại cung đối Điền Trạch",
    "Thiên Tướng đồng cung Thất Sát tại cung đối Điền Trạch",
    "Thiên Tướng đồng cung Phá Quân tại cung đối Điền Trạch",
    "Thiên Lương đồng cung Thất Sát tại cung đối Điền Trạch",
    "Thiên Lương đồng cung Phá Quân tại cung đối Điền Trạch",
    "Thất Sát đồng cung Phá Quân tại cung đối Điền Trạch"
)

$startRow = 4305
$count = $newStrings.Count

$arr = New-Object 'object[,]' $count,2
for ($i = 0; $i -lt $count; $i++) {
    $arr[$i,0] = $newStrings[$i]
    $arr[$i,1] = $newStrings[$i]
}

$endRow = $startRow + $count - 1
$rng = $ws.Range("A$startRow" + ":B$endRow")
$rng.Value = $arr

# Match the author's final selection/viewport on the newly added block.
$ws.Activate()
$null = $ws.Range("B4319:B4409").Select()
